# Apply updated GroupGrowthRateMax values: freeze 5-year growth rate at its
# peak value (the rate in effect just before the first drop in the logistic
# diffusion projection), and update the explanatory note text to reflect
# this methodology change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GroupGrowthRateMax")

# New "rate" (column D) values for rows 2-145, keyed by row number.
$newRates = @{
    2 = 0.965261432;
    3 = 1.291488707;
    4 = 1.291488707;
    5 = 1.291488707;
    6 = 1.291488707;
    7 = 1.291488707;
    8 = 0.97267241;
    9 = 1.314564727;
    10 = 1.314564727;
    11 = 1.314564727;
    12 = 1.314564727;
    13 = 1.314564727;
    14 = 0.973374088;
    15 = 1.316770982;
    16 = 1.316770982;
    17 = 1.316770982;
    18 = 1.316770982;
    19 = 1.316770982;
    20 = 0.980715628;
    21 = 1.340080203;
    22 = 1.340080203;
    23 = 1.340080203;
    24 = 1.340080203;
    25 = 1.340080203;
    26 = 0.985348725;
    27 = 1.355004872;
    28 = 1.355004872;
    29 = 1.355004872;
    30 = 1.355004872;
    31 = 1.355004872;
    32 = 0.985348725;
    33 = 1.355004872;
    34 = 1.355004872;
    35 = 1.355004872;
    36 = 1.355004872;
    37 = 1.355004872;
    38 = 0.979617728;
    39 = 1.336568026;
    40 = 1.336568026;
    41 = 1.336568026;
    42 = 1.336568026;
    43 = 1.336568026;
    44 = 0.985225304;
    45 = 1.354605115;
    46 = 1.354605115;
    47 = 1.354605115;
    48 = 1.354605115;
    49 = 1.354605115;
    50 = 0.98545597;
    51 = 1.355352333;
    52 = 1.355352333;
    53 = 1.355352333;
    54 = 1.355352333;
    55 = 1.355352333;
    56 = 0.98545597;
    57 = 1.355352333;
    58 = 1.355352333;
    59 = 1.355352333;
    60 = 1.355352333;
    61 = 1.355352333;
    62 = 0.98545597;
    63 = 1.355352333;
    64 = 1.355352333;
    65 = 1.355352333;
    66 = 1.355352333;
    67 = 1.355352333;
    68 = 0.985105904;
    69 = 1.354218491;
    70 = 1.354218491;
    71 = 1.354218491;
    72 = 1.354218491;
    73 = 1.354218491;
    74 = 0.985105904;
    75 = 1.354218491;
    76 = 1.354218491;
    77 = 1.354218491;
    78 = 1.354218491;
    79 = 1.354218491;
    80 = 0.985105904;
    81 = 1.354218491;
    82 = 1.354218491;
    83 = 1.354218491;
    84 = 1.354218491;
    85 = 1.354218491;
    86 = 0.985105904;
    87 = 1.354218491;
    88 = 1.354218491;
    89 = 1.354218491;
    90 = 1.354218491;
    91 = 1.354218491;
    92 = 0.98527064;
    93 = 1.354751936;
    94 = 1.354751936;
    95 = 1.354751936;
    96 = 1.354751936;
    97 = 1.354751936;
    98 = 0.98527064;
    99 = 1.354751936;
    100 = 1.354751936;
    101 = 1.354751936;
    102 = 1.354751936;
    103 = 1.354751936;
    104 = 0.98527064;
    105 = 1.354751936;
    106 = 1.354751936;
    107 = 1.354751936;
    108 = 1.354751936;
    109 = 1.354751936;
    110 = 0.98527064;
    111 = 1.354751936;
    112 = 1.354751936;
    113 = 1.354751936;
    114 = 1.354751936;
    115 = 1.354751936;
    116 = 0.971043362;
    117 = 1.309456939;
    118 = 1.309456939;
    119 = 1.309456939;
    120 = 1.309456939;
    121 = 1.309456939;
    122 = 0.984146192;
    123 = 1.35111499;
    124 = 1.35111499;
    125 = 1.35111499;
    126 = 1.35111499;
    127 = 1.35111499;
    128 = 0.985381484;
    129 = 1.355110999;
    130 = 1.355110999;
    131 = 1.355110999;
    132 = 1.355110999;
    133 = 1.355110999;
    134 = 0.985450374;
    135 = 1.355334198;
    136 = 1.355334198;
    137 = 1.355334198;
    138 = 1.355334198;
    139 = 1.355334198;
    140 = 0.985450374;
    141 = 1.355334198;
    142 = 1.355334198;
    143 = 1.355334198;
    144 = 1.355334198;
    145 = 1.355334198
}

$oldPhrase = "logistic diffusion model;"
$newPhrase = "logistic diffusion model, 5-year growth is frozen before first drop in projections;"

foreach ($row in $newRates.Keys) {
    # Update the growth-rate value in column D
    $ws.Cells.Item($row, 4).Value = $newRates[$row]

    # Update the explanatory note in column E to mention the new freezing methodology
    $noteCell = $ws.Cells.Item($row, 5)
    $note = $noteCell.Text
    if ($note -ne $null -and $note.Contains($oldPhrase) -and -not $note.Contains($newPhrase)) {
        $noteCell.Value = $note.Replace($oldPhrase, $newPhrase)
    }
}
